# Add agenda and purpose slides.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Clear the leftover (mismatched) speaker notes text on slide 1 --
#    the notes placeholder keeps an empty paragraph afterwards.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$notesShape = $s1.NotesPage.Shapes.Item(2)
$notesShape.TextFrame.TextRange.Text = ""

# ---------------------------------------------------------------------
# 2. New slide 3 - "Agenda" (Title and Content layout)
# ---------------------------------------------------------------------
$sAgenda = $p.Slides.Add(3, 2)
$sAgenda.Shapes.Item(1).TextFrame.TextRange.Text = "Agenda"

$trAgenda = $sAgenda.Shapes.Item(2).TextFrame.TextRange
$trAgenda.Text = "Purpose of the transaction log"
[void]$trAgenda.InsertAfter("`rOrganization of the transaction log")
[void]$trAgenda.InsertAfter("`rFlushing the log")
[void]$trAgenda.InsertAfter("`rClearing the log")
[void]$trAgenda.InsertAfter("`rCheckpoints")
[void]$trAgenda.InsertAfter("`rRollback operations")
[void]$trAgenda.InsertAfter("`rVLF fragmentation")
[void]$trAgenda.InsertAfter("`rLog monitoring")

# ---------------------------------------------------------------------
# 3. New slide 4 - "Purpose of the Transaction Log" (Title and Content layout)
# ---------------------------------------------------------------------
$sPurpose = $p.Slides.Add(4, 2)
$sPurpose.Shapes.Item(1).TextFrame.TextRange.Text = "Purpose of the Transaction Log"

$bodyPurpose = $sPurpose.Shapes.Item(2)
$bodyPurpose.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

$trPurpose = $bodyPurpose.TextFrame.TextRange
$trPurpose.Text = "Primary purposes"
[void]$trPurpose.InsertAfter("`rDurability")
[void]$trPurpose.InsertAfter("`rCrash recovery / restore operations")
[void]$trPurpose.InsertAfter("`rAtomicity")
[void]$trPurpose.InsertAfter("`rThought experiment")
[void]$trPurpose.InsertAfter("`rWhat would SQL be like without ")
[void]$trPurpose.InsertAfter("a transaction log?")
[void]$trPurpose.InsertAfter("`rSecondary purposes")
[void]$trPurpose.InsertAfter("`rLog reader (replication, CDC)")
[void]$trPurpose.InsertAfter("`rMirroring / Availability Groups / log shipping")
[void]$trPurpose.InsertAfter("`rSnapshots")

$trPurpose.Paragraphs(2,1).IndentLevel = 2   # Durability
$trPurpose.Paragraphs(3,1).IndentLevel = 3   # Crash recovery / restore operations
$trPurpose.Paragraphs(4,1).IndentLevel = 2   # Atomicity
$trPurpose.Paragraphs(5,1).IndentLevel = 2   # Thought experiment
$trPurpose.Paragraphs(6,1).IndentLevel = 3   # What would SQL be like without a transaction log?
$trPurpose.Paragraphs(8,1).IndentLevel = 2   # Log reader (replication, CDC)
$trPurpose.Paragraphs(9,1).IndentLevel = 2   # Mirroring / Availability Groups / log shipping
$trPurpose.Paragraphs(10,1).IndentLevel = 2  # Snapshots
